# Daily attendance processing - 2026-01-25 19:10:07
# Reorders the "Recorded By" (column G) contributor lists so that
# "System" is listed before external/backup email addresses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "system, System, backup@backdoor.com"
    3   = "System, dnasr281@gmail.com"
    4   = "System, backup@backdoor.com"
    5   = "System, backup@backdoor.com"
    6   = "System, dnasr281@gmail.com"
    8   = "System, backup@backdoor.com"
    28  = "system, System, backup@backdoor.com"
    29  = "System, dnasr281@gmail.com"
    30  = "System, backup@backdoor.com"
    31  = "System, backup@backdoor.com"
    32  = "System, dnasr281@gmail.com"
    34  = "System, backup@backdoor.com"
    54  = "system, System, backup@backdoor.com"
    55  = "System, dnasr281@gmail.com"
    56  = "System, backup@backdoor.com"
    57  = "System, backup@backdoor.com"
    58  = "System, dnasr281@gmail.com"
    60  = "System, backup@backdoor.com"
    80  = "System, backup@backdoor.com"
    81  = "System, backup@backdoor.com"
    82  = "System, backup@backdoor.com"
    87  = "admin@admin.com, dnasr281@gmail.com"
    106 = "System, backup@backdoor.com"
    107 = "System, backup@backdoor.com"
    108 = "System, backup@backdoor.com"
    113 = "admin@admin.com, dnasr281@gmail.com"
    132 = "System, backup@backdoor.com"
    133 = "System, backup@backdoor.com"
    134 = "System, backup@backdoor.com"
    139 = "admin@admin.com, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
